# Generate Report for Handoff
#
# The file "89748552-4df1-4d18-88ae-347e9f84a6ea.md" finished its zh-cn and
# de-de localization handoff and is now "Ready for handoff". Update the
# per-language sheets (zh-cn, de-de) and the roll-up "Overview" sheet with
# the new status + timestamps, and widen the Status columns to fit the new
# (longer) status text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the 89748552 file.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("H3").Value = "2016-08-26 18:14:43"

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the 89748552 file.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("H3").Value = "2016-08-26 18:14:47"

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the 89748552 file. zh-cn status (E), de-de
# status (F), and the latest handoff-xliff-generate date (G) all roll
# up to "Ready for handoff" / the newest timestamp.
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("E3").Value = "Ready for handoff"
$wsOv.Range("F3").Value = "Ready for handoff"
$wsOv.Range("G3").Value = "2016-08-26 18:14:47"

# ---------------------------------------------------------------------
# The Status column on each sheet grew wider to fit "Ready for handoff".
# ---------------------------------------------------------------------
$wsOv.Columns.Item(5).ColumnWidth = 16.3
$wsOv.Columns.Item(6).ColumnWidth = 16.3
$wsZh.Columns.Item(3).ColumnWidth = 16.3
$wsDe.Columns.Item(3).ColumnWidth = 16.3
